# Update the "人数" (F column) figures on both the "展览" and "全部类型"
# worksheets to the new values captured at the later scrape (456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value (identical update on both sheets)
$updates = @{
    3  = 51
    5  = 91
    7  = 1268
    8  = 1537
    9  = 341
    10 = 394
    12 = 152
    14 = 63
    15 = 107
    17 = 306
    18 = 324
    19 = 1735
    23 = 669
    25 = 336
    26 = 4187
    29 = 1090
    32 = 548
    34 = 253
    36 = 139
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
